$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = "2x DXW D4250 600KV 3-7S Outrunner"
$ws.Range("B25").Value = "2 x ESC (Bộ điều tốc) 100A"
$ws.Range("B29").Value = "2x Lipo CNHL LVNCell 6s 5200mah 65C"
$ws.Range("C29").Value = "Nguồn chính Cấu hình 6S2P - 10400mAh"
$ws.Range("B36").Value = "Thanh carbon phi 16mm"
$ws.Range("B37").Value = "Nhựa in 3D"

$ws.Range("C41").Select()
